$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.436.48"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "3.686.41"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'686.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "

$ws.Range("D6").Value = "'159.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.00%  "

$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").Value = "'0.493"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  -1.52%  "

$ws.Range("E10").Value = "  -2.17%  "

$ws.Range("D11").Value = "'0.434"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.52%  "

$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").Value = "4.310.32"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "'32.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.38%  "

$ws.Range("D15").Value = "3.709.86"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("D16").Value = "69.426.60"
$ws.Range("E16").Value = "  +0.11%  "

$ws.Range("E17").Value = "  +1.75%  "

$ws.Range("D18").Value = "'15.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.49%  "

$ws.Range("D19").Value = "'6.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.97%  "

$ws.Range("D20").Value = "'471.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").Value = "'9.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.06%  "

$ws.Range("E22").Value = "  -2.11%  "

$ws.Range("D23").Value = "'79.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "

$ws.Range("D24").Value = "3.835.58"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").Value = "'0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.88%  "

$ws.Range("D27").Value = "'10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.10%  "

$ws.Range("D28").Value = "'9.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.94%  "

$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("D30").Value = "'1.74"
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = "  -5.53%  "

$ws.Range("D32").Value = "'6.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'0.998"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.61%  "

$ws.Range("D35").Value = "3.663.36"
$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("D36").Value = "'0.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").Value = "'8.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.90%  "

$ws.Range("D38").Value = "'6.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "

$ws.Range("E40").Value = "  +1.81%  "

$ws.Range("D41").Value = "'0.0897"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.92%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").Value = "'0.942"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "

$ws.Range("D44").Value = "'165.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.53%  "

$ws.Range("D45").Value = "'47.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("D46").Value = "'0.000283"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("E47").Value = "  -3.76%  "

$ws.Range("D48").Value = "'1.11"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.13%  "

$ws.Range("E49").Value = "  +0.41%  "

$ws.Range("D50").Value = "'27.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.32%  "

$ws.Range("D51").Value = "'7.76"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.57%  "
